# Applies the gh-pages data refresh (commit 456a3b4) to both the
# "展览" (Exhibitions) and "全部类型" (All types) sheets.
#
# Summary of changes (identical pattern on both sheets):
#  1. "南宁·星STAR国潮嘉年华" is marked cancelled and its lowest price
#     becomes the text "不可售" instead of a numeric 50.
#  2. "南宁·小蜜蜂动漫嘉年华2.0" interest count 263 -> 264.
#  3. A brand new event row is inserted right before "南宁·0713国乙ONLY":
#     "南宁·漫控嘉年华09暨南宁高校动漫联盟十六周年联合漫展" on 2024-07-12.
#  4. "南宁·AB动漫游戏嘉年华" interest count 1941 -> 1945.
#  5. "南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）" interest count 4550 -> 4558.
#  6. "南宁·火影忍者only" interest count 71 -> 72.
#  7. "南宁·蔚蓝档案only" interest count 319 -> 321.
#
# NOTE: this runtime's PowerShell engine does not bind named (-Param)
# arguments correctly, so every helper function below uses purely
# positional parameters.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

function Insert-NewEventRow {
    param($ws, $RowIndex, $Date, $Title, $Location, $TimeRange, $WantCount, $MinPrice, $Link, $Cover)

    # Push everything from RowIndex downward by one row.
    $ws.Rows.Item($RowIndex).Insert()

    # Copy the number/border formatting of column A from the row that is
    # now directly below (it still carries the original "index" style),
    # so the new row's A cell matches the look of every other data row.
    $belowRow = $RowIndex + 1
    $ws.Range("A" + $belowRow).Copy()
    $ws.Range("A" + $RowIndex).PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false

    # The "index" column (A) is the running number = row above's index + 1.
    $aboveIndexCell = $ws.Range("A" + ($RowIndex - 1))
    $newIndex = [int]$aboveIndexCell.Value2 + 1

    $ws.Range("A" + $RowIndex).Value = $newIndex
    $ws.Range("B" + $RowIndex).Value = $Date
    $ws.Range("C" + $RowIndex).Value = $Title
    $ws.Range("D" + $RowIndex).Value = $Location
    $ws.Range("E" + $RowIndex).Value = $TimeRange
    $ws.Range("F" + $RowIndex).Value = $WantCount
    $ws.Range("G" + $RowIndex).Value = $MinPrice
    $ws.Range("H" + $RowIndex).Value = $Link
    $ws.Range("I" + $RowIndex).Value = $Cover
}

function Update-Sheet {
    param($ws, $StarRow, $BeeRow, $NewEventRow, $AbRow, $LiangyaRow, $NarutoRow, $AzurRow)

    # 1. 南宁·星STAR国潮嘉年华 -> mark as cancelled, price becomes "不可售".
    $ws.Range("C" + $StarRow).Value = "南宁·星STAR国潮嘉年华（取消）"
    $ws.Range("G" + $StarRow).Value = "不可售"

    # 2. 南宁·小蜜蜂动漫嘉年华2.0 想去人数 263 -> 264.
    $ws.Range("F" + $BeeRow).Value = 264

    # 3. Insert the new "漫控嘉年华09" event before "0713国乙ONLY".
    Insert-NewEventRow $ws $NewEventRow "2024-07-12" "南宁·漫控嘉年华09暨南宁高校动漫联盟十六周年联合漫展" "民族大道106号 南宁国际会展中心" "2024.07.12 09:30-07.14 17:00" 8 50 "https://show.bilibili.com/platform/detail.html?id=87182" "//i1.hdslb.com/bfs/openplatform/202406/x4UZPn301718159475475.jpeg"

    # Row indices below the inserted row have all shifted down by one.
    $AbRow = $AbRow + 1
    $LiangyaRow = $LiangyaRow + 1
    $NarutoRow = $NarutoRow + 1
    $AzurRow = $AzurRow + 1

    # 4. 南宁·AB动漫游戏嘉年华 想去人数 1941 -> 1945.
    $ws.Range("F" + $AbRow).Value = 1945

    # 5. 南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典） 想去人数 4550 -> 4558.
    $ws.Range("F" + $LiangyaRow).Value = 4558

    # 6. 南宁·火影忍者only 想去人数 71 -> 72.
    $ws.Range("F" + $NarutoRow).Value = 72

    # 7. 南宁·蔚蓝档案only 想去人数 319 -> 321.
    $ws.Range("F" + $AzurRow).Value = 321
}

# ---- Sheet "展览" ----
$wsExpo = $wb.Worksheets.Item("展览")
Update-Sheet $wsExpo 2 4 6 8 10 11 12

# ---- Sheet "全部类型" ----
$wsAll = $wb.Worksheets.Item("全部类型")
Update-Sheet $wsAll 3 6 8 12 14 15 16
